$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.254.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.24%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.171.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -8.33%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'564.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.93%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'169.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.94%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'0.606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.169.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -8.33%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -6.38%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -5.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.721.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -8.32%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.83%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'27.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -9.25%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.233.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.11%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -5.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.174.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -8.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.83%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.22%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'352.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.25%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -5.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'68.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -5.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -4.67%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -5.80%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -3.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.13%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.25%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'5.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -5.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -5.39%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'22.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -7.17%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Fetch.AI"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.91%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Aptos"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'6.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.93%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -6.88%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'154.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.21%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.818"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'25.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.98%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.99%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.616.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.56%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -7.22%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -6.74%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'39.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0654"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.89%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'23.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'321.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.03%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -7.84%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.56%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.02%  "
$ws.Range("E51").Style = "Normal"
